$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second row of the table (the "H 72" record that had a blank
# final column), shifting all subsequent data rows up by one.
$ws.Rows.Item(2).Delete()
